# Update gh-pages to output generated at 456a3b4
# Applies the numeric "want-to-go" / price bumps across the four sheets, and
# inserts the new "「多厨狂喜」白金交响乐团二次元交响音乐会" (2024-11-03) row into
# the "全部类型" aggregate sheet (pushing the later rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (F-column "want to go" counters)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2443
$ws1.Range("F7").Value = 295
$ws1.Range("F9").Value = 3407
$ws1.Range("F14").Value = 15
$ws1.Range("F16").Value = 998
$ws1.Range("F17").Value = 1748
$ws1.Range("F19").Value = 414
$ws1.Range("F20").Value = 1523
$ws1.Range("F21").Value = 1063
$ws1.Range("F22").Value = 95
$ws1.Range("F23").Value = 130
$ws1.Range("F24").Value = 4090
$ws1.Range("F27").Value = 1186

# ---------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G3").Value = 220
$ws2.Range("F23").Value = 10
$ws2.Range("F26").Value = 62
$ws2.Range("F36").Value = 58
$ws2.Range("F39").Value = 380
$ws2.Range("F47").Value = 15
$ws2.Range("F48").Value = 15

# ---------------------------------------------------------------
# Sheet "本地生活"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F12").Value = 2880
$ws3.Range("F13").Value = 416
$ws3.Range("F14").Value = 742
$ws3.Range("F15").Value = 121

# ---------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 416
$ws4.Range("F9").Value = 742
$ws4.Range("G10").Value = 220
$ws4.Range("F11").Value = 121
$ws4.Range("F15").Value = 295
$ws4.Range("F17").Value = 3407
$ws4.Range("F23").Value = 15
$ws4.Range("F25").Value = 998
$ws4.Range("F28").Value = 10
$ws4.Range("F29").Value = 1748
$ws4.Range("F30").Value = 414
$ws4.Range("F32").Value = 1523
# this event's own "want to go" count also bumps 59 -> 62
$ws4.Range("F34").Value = 62

# Insert a new row 35 duplicating the (now-updated) row 34 event -- same
# listing re-surfaced for its 2024-11-03 show date -- then push every row
# from the old 35 ("2024年刘明月专场生日会") through 50 down by one, so the
# sheet grows from 50 to 51 data rows.
$ws4.Rows.Item(35).Insert()

# match formatting (bold/centered index style + borders) of the index column
$ws4.Range("A34").Copy()
$ws4.Range("A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Range("A35").Value = 34
# force text (not auto-converted to a date serial) for the date-looking label,
# then restore the original "General" number format
$ws4.Range("B35").NumberFormat = "@"
$ws4.Range("B35").Value = "2024-11-03"
$ws4.Range("B35").NumberFormat = "General"
$ws4.Range("C35").Value = "上海·「多厨狂喜」白金交响乐团二次元交响音乐会"
$ws4.Range("D35").Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$ws4.Range("E35").Value = "2024.11.03 14:00-11.03 16:00"
$ws4.Range("F35").Value = 62
$ws4.Range("G35").Value = 99
$ws4.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=93086"
$ws4.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202409/GvAAxiwb1727619935967.jpeg"
